$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the date value for A13 (2012-08-08, serial 41129) - keep existing style (format already applied)
$ws.Range("A13").Value = (Get-Date -Year 2012 -Month 8 -Day 8 -Hour 0 -Minute 0 -Second 0).Date

# Set the activity text for B13 (new shared string entry)
$ws.Range("B13").Value = "Tried to get OpenSSL custom OpenCL engine running"

# Update the active cell / selection to B14, matching the author's cursor move after data entry
$ws.Activate()
$ws.Range("B14").Select()
